$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.115.59'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '4.020.09'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.28'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.28'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.705'
$ws.Range('E7').Value = '  +12.52%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.750'
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000327'
$ws.Range('E11').Value = '  -4.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.52'
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.64'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.660.52'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.036.40'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.08'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.53'
$ws.Range('E17').Value = '  -4.76%  '
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.052.34'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '428.64'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '98.24'
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.25'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.39'
$ws.Range('E26').Value = '  -8.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.72'
$ws.Range('E27').Value = '  -5.14%  '
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.77'
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.61'
$ws.Range('E30').Value = '  +17.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.39'
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.20'
$ws.Range('E33').Value = '  +3.64%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '683.13'
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '45.05'
$ws.Range('E35').Value = '  +9.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.64'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.449'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0828'
$ws.Range('E38').Value = '  -9.32%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.150'
$ws.Range('E39').Value = '  -3.84%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.43'
$ws.Range('E40').Value = '  -5.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.27'
$ws.Range('E43').Value = '  +5.27%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0487'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.72'
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.92'
$ws.Range('E47').Value = '  +7.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.40'
$ws.Range('E48').Value = '  -3.60%  '
$ws.Range('E49').Value = '  -5.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000269'
$ws.Range('E50').Value = '  -4.43%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.25'
$ws.Range('E51').Value = '  -2.95%  '
